$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4527.612
$ws.Range("I15").Value = 4527.612
$ws.Range("K15").Value = 13582.836
$ws.Range("M15").Value = -13413.836
$ws.Range("H17").Value = 12198722
$ws.Range("J17").Value = 12824253
$ws.Range("L17").Value = 38472759
$ws.Range("N17").Value = -38473095
$ws.Range("H51").Value = 4117340.5
$ws.Range("I51").Value = 55556456
$ws.Range("J51").Value = 2211.2
$ws.Range("K51").Value = 55556456
$ws.Range("L51").Value = 2211.2
$ws.Range("M51").Value = -55555972
$ws.Range("N51").Value = -3179.2
$ws.Range("H74").Value = 3222.9412
$ws.Range("I74").Value = 3445.4546
$ws.Range("J74").Value = 2815
$ws.Range("K74").Value = 3445.4546
$ws.Range("L74").Value = 2815
$ws.Range("M74").Value = -2509.4546
$ws.Range("N74").Value = -4687
$ws.Range("H77").Value = 3222.9412
$ws.Range("I77").Value = 3445.4546
$ws.Range("J77").Value = 2815
$ws.Range("K77").Value = 17227.273
$ws.Range("L77").Value = 14075
$ws.Range("M77").Value = -12547.273
$ws.Range("N77").Value = -23435
$ws.Range("H92").Value = 908
$ws.Range("I92").Value = 908
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 908
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 340
$ws.Range("N92").ClearContents()
$ws.Range("H94").Value = 4304.5386
$ws.Range("I94").Value = 4304.5386
$ws.Range("K94").Value = 4304.5386
$ws.Range("M94").Value = -3853.5386
$ws.Range("H96").Value = 66674468
$ws.Range("I96").Value = 8808.25
$ws.Range("J96").Value = 90916530
$ws.Range("K96").Value = 26424.75
$ws.Range("L96").Value = 272749590
$ws.Range("M96").Value = -25051.75
$ws.Range("N96").Value = -272752336
$ws.Range("H100").Value = 2329.2856
$ws.Range("I100").Value = 2329.2856
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2329.2856
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1788.2856
$ws.Range("N100").ClearContents()
$ws.Range("H101").Value = 898.3333
$ws.Range("J101").Value = 1207.8572
$ws.Range("L101").Value = 3623.5716
$ws.Range("N101").Value = -6867.571599999999
$ws.Range("H129").Value = 2558.087
$ws.Range("I129").Value = 3731.3333
$ws.Range("J129").Value = 2382.1
$ws.Range("K129").Value = 11193.9999
$ws.Range("L129").Value = 7146.299999999999
$ws.Range("M129").Value = -6193.999899999999
$ws.Range("N129").Value = -17146.3
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 41668264
$ws.Range("I45").Value = 55556990
$ws.Range("J45").Value = 2083.8333
$ws.Range("K45").Value = 55556990
$ws.Range("L45").Value = 2083.8333
$ws.Range("M45").Value = -55556613
$ws.Range("N45").Value = -2837.8333
$ws.Range("H96").Value = 25099.666
$ws.Range("J96").Value = 25099.666
$ws.Range("L96").Value = 25099.666
$ws.Range("N96").Value = -30591.666
$ws.Range("H97").Value = 779.79486
$ws.Range("I97").Value = 594.1111
$ws.Range("J97").Value = 1197.5834
$ws.Range("K97").Value = 594.1111
$ws.Range("L97").Value = 1197.5834
$ws.Range("M97").Value = -98.11109999999996
$ws.Range("N97").Value = -2189.5834
$ws.Range("H102").Value = 12987.81
$ws.Range("I102").Value = 2398.4443
$ws.Range("J102").Value = 20929.834
$ws.Range("K102").Value = 2398.4443
$ws.Range("L102").Value = 20929.834
$ws.Range("M102").Value = -776.4443000000001
$ws.Range("N102").Value = -24173.834
$ws.Range("H119").Value = 52619
$ws.Range("J119").Value = 52619
$ws.Range("L119").Value = 52619
$ws.Range("N119").Value = -62295
$ws.Range("H122").Value = 1968.6857
$ws.Range("I122").Value = 1998.4814
$ws.Range("J122").Value = 1868.125
$ws.Range("K122").Value = 5995.4442
$ws.Range("L122").Value = 5604.375
$ws.Range("M122").Value = -3545.4442
$ws.Range("N122").Value = -10504.375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 580.52
$ws.Range("I94").Value = 587.0909
$ws.Range("K94").Value = 587.0909
$ws.Range("M94").Value = -136.0909
$ws.Range("H99").Value = 2217.2778
$ws.Range("I99").Value = 2170.5881
$ws.Range("K99").Value = 2170.5881
$ws.Range("M99").Value = -672.5880999999999
$ws.Range("H105").Value = 2729.5
$ws.Range("I105").Value = 2441.3333
$ws.Range("J105").Value = 4170.3335
$ws.Range("K105").Value = 2441.3333
$ws.Range("L105").Value = 4170.3335
$ws.Range("M105").Value = -694.3332999999998
$ws.Range("N105").Value = -7664.3335
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3358.8
$ws.Range("I105").Value = 3565.4443
$ws.Range("J105").Value = 1499
$ws.Range("K105").Value = 3565.4443
$ws.Range("L105").Value = 1499
$ws.Range("M105").Value = -1818.4443
$ws.Range("N105").Value = -4993
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4646.5454
$ws.Range("I3").Value = 1924.6666
$ws.Range("J3").Value = 7912.8
$ws.Range("K3").Value = 5773.9998
$ws.Range("L3").Value = 23738.4
$ws.Range("M3").Value = -5661.9998
$ws.Range("N3").Value = -23962.4
$ws.Range("H109").Value = 2461.8
$ws.Range("I109").Value = 860.4
$ws.Range("J109").Value = 3262.5
$ws.Range("K109").Value = 2581.2
$ws.Range("L109").Value = 9787.5
$ws.Range("M109").Value = -1541.2
$ws.Range("N109").Value = -11867.5
$ws.Range("H137").Value = 35724660
$ws.Range("I137").Value = 4464
$ws.Range("J137").Value = 43489920
$ws.Range("K137").Value = 13392
$ws.Range("L137").Value = 130469760
$ws.Range("M137").Value = -8292
$ws.Range("N137").Value = -130479960
$ws.Range("H141").Value = 43481944
$ws.Range("I141").Value = 71432060
$ws.Range("J141").Value = 3987.7778
$ws.Range("K141").Value = 214296180
$ws.Range("L141").Value = 11963.3334
$ws.Range("M141").Value = -214291000
$ws.Range("N141").Value = -22323.3334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 34485332
$ws.Range("I7").Value = 43480350
$ws.Range("J7").Value = 4417.5
$ws.Range("K7").Value = 43480350
$ws.Range("L7").Value = 4417.5
$ws.Range("M7").Value = -43480238
$ws.Range("N7").Value = -4641.5
$ws.Range("H40").Value = 2519.8333
$ws.Range("I40").Value = 2005.3529
$ws.Range("K40").Value = 2005.3529
$ws.Range("M40").Value = -1869.3529
$ws.Range("H81").Value = 42181
$ws.Range("J81").Value = 42181
$ws.Range("L81").Value = 42181
$ws.Range("N81").Value = -44177
$ws.Range("H84").Value = 42181
$ws.Range("J84").Value = 42181
$ws.Range("L84").Value = 126543
$ws.Range("N84").Value = -136527
$ws.Range("H93").Value = 1479.0834
$ws.Range("I93").Value = 1398.25
$ws.Range("J93").Value = 1519.5
$ws.Range("K93").Value = 1398.25
$ws.Range("L93").Value = 1519.5
$ws.Range("M93").Value = -150.25
$ws.Range("N93").Value = -4015.5
$ws.Range("H126").Value = 34485332
$ws.Range("I126").Value = 43480350
$ws.Range("J126").Value = 4417.5
$ws.Range("K126").Value = 130441050
$ws.Range("L126").Value = 13252.5
$ws.Range("M126").Value = -130438580
$ws.Range("N126").Value = -18192.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2395
$ws.Range("I96").Value = 999
$ws.Range("J96").Value = 3325.6667
$ws.Range("K96").Value = 3325.6667
$ws.Range("L96").Value = 3325.6667
$ws.Range("M96").Value = 374
$ws.Range("N96").Value = -6071.6667
